# Generate Report for Handback
# Updates the "Correspond Handoff/Handback Datetime" and
# "Latest HO Xliff Generate Date" timestamps on the handback-status report.

$wb = $excel.ActiveWorkbook

# "Overview" sheet: G2 "Latest HO Xliff Generate Date" for the first file
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-22 16:46:59"

# "zh-cn" sheet: H2 "Correspond Handoff Datetime" and K2 "Correspond Handback DateTime"
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-22 16:46:55"
$wsZhCn.Range("K2").Value = "2016-08-22 16:47:21"

# "de-de" sheet: H2 shares the same source timestamp as Overview!G2 (both
# reflect the handoff time for 6c7624f3-...md), and K2 is its own
# "Correspond Handback DateTime".
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-22 16:46:59"
$wsDeDe.Range("K2").Value = "2016-08-22 16:47:28"
